# Stage 1: update companies data
#
# The underlying company records (Company Name, Company Number, Category,
# SIC Codes, SIC Description, Typical Use Case -- i.e. columns A, B, H, I,
# J, K) have been re-ordered across rows 3-11 (row 8 / "THE DISLEY GROUP
# LTD" is unaffected). The per-row discovery metadata in columns C, D, E,
# F, G stays put on each row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New company-record values keyed by destination row number.
$data = @{
    3  = @{ A = "GANDER INVESTMENTS LTD";                  B = "16473515"; H = "Investments"; I = "68100,68209";       J = "";                                                  K = "" }
    4  = @{ A = "SEVEN (HOLDCO) LIMITED";                   B = "16473606"; H = "Other";       I = "64209";            J = "Activities of other holding companies n.e.c.";      K = "Catch-all SPV: protected cells, cell companies, bespoke feeder vehicles." }
    5  = @{ A = "BRIDGEWICK PARTNERS LIMITED";              B = "16473142"; H = "Partners";     I = "64999";            J = "Financial intermediation not elsewhere classified"; K = "Catch-all credit-oriented SPVs for novel lending structures." }
    6  = @{ A = "MARMIMI HOLDING LIMITED";                  B = "16473234"; H = "Other";        I = "64209";            J = "Activities of other holding companies n.e.c.";      K = "Catch-all SPV: protected cells, cell companies, bespoke feeder vehicles." }
    7  = @{ A = "AJ INVESTMENT AND CONSULTANCY LTD";        B = "16473328"; H = "Investments";  I = "64306,70229";      J = "Activities of real estate investment trusts";       K = "UK-regulated REIT companies." }
    9  = @{ A = "TLJ INVESTMENT LTD";                       B = "16473151"; H = "Investments";  I = "41100,55100,68100"; J = "";                                                  K = "" }
    10 = @{ A = "INTERCONTINENTAL HOLDING COMPANY LIMITED"; B = "16473418"; H = "Other";        I = "64209";            J = "Activities of other holding companies n.e.c.";      K = "Catch-all SPV: protected cells, cell companies, bespoke feeder vehicles." }
    11 = @{ A = "GAUNT CAPITAL LTD";                        B = "16473262"; H = "Capital";       I = "64209";            J = "Activities of other holding companies n.e.c.";      K = "Catch-all SPV: protected cells, cell companies, bespoke feeder vehicles." }
}

# Columns that may contain purely-numeric-looking text (company numbers,
# single-code SIC codes) which Excel would otherwise auto-coerce into a
# numeric cell. Force them to remain text, then reset the style back to
# Normal so we don't leave a stray number-format style on the cell.
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

foreach ($row in $data.Keys) {
    $rec = $data[$row]

    Set-TextValue $ws.Range("A$row") $rec.A
    Set-TextValue $ws.Range("B$row") $rec.B
    $ws.Range("H$row").Value = $rec.H
    Set-TextValue $ws.Range("I$row") $rec.I
    $ws.Range("J$row").Value = $rec.J
    $ws.Range("K$row").Value = $rec.K
}
